$wb = $excel.ActiveWorkbook

# ----- Sheet "配置" (first sheet) -----
$ws1 = $wb.Worksheets.Item(1)

# Remove the last row (the sheet shrinks from 8 to 7 rows)
$ws1.Rows.Item(8).Delete()

# Update remaining rows with the new Ollama-related configuration
$ws1.Range("B1").Value = "内容"

$ws1.Range("A2").Value = "启用AI"
$ws1.Range("B2").Value = "Y"

$ws1.Range("A3").Value = "Ollama Host"
$ws1.Range("B3").Value = "192.168.96.246"

$ws1.Range("A4").Value = "Ollama Port"
$ws1.Range("B4").Value = 11434

$ws1.Range("A5").Value = "Ollama Model"
$ws1.Range("B5").Value = "wangshenzhi/llama3-8b-chinese-chat-ollama-q8"

$ws1.Range("A6").Value = "利用语义分析删除相关内容"
$ws1.Range("B6").Value = "Y"

$ws1.Range("A7").Value = "生成AI处理前后对比文档"
$ws1.Range("B7").Value = "Y"

$ws1.Range("B5").Select()

# Page setup now targets A4-equivalent paper in portrait orientation
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ----- Sheet "语义分析" (second sheet) -----
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "表情符号"
$ws2.Range("B2").Value = 8

$ws2.Range("A3").Value = "投票"
$ws2.Range("B3").Value = 8

# The second sheet becomes the active / selected tab
$ws2.Activate()
$ws2.Range("D12").Select()
